$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPDS-JP")

# Fill in card names for two existing rows that were missing column A
$ws.Range("A33").Value = "Darklord Mastema"
$ws.Range("A36").Value = "Altar of the Darklords"

# Append 4 new rows (38-41) for the new cards/ids, following the existing
# B/C/E pattern used throughout the sheet.
$ws.Range("A38").Value = "Darklords' Temptation"
$ws.Range("B38").Value = 100405037
$ws.Range("C38").Value = ":"
$ws.Range("E38").Value = ","

$ws.Range("B39").Value = 100405038
$ws.Range("C39").Value = ":"
$ws.Range("E39").Value = ","

$ws.Range("B40").Value = 100405039
$ws.Range("C40").Value = ":"
$ws.Range("E40").Value = ","

$ws.Range("B41").Value = 100405040
$ws.Range("C41").Value = ":"
$ws.Range("E41").Value = ","

# Match the final cursor/selection position from the source workbook.
$ws.Range("F45").Select() | Out-Null
